# Update "想去人数" (want-to-go count) figures in column F
# for both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1024
    3  = 288
    4  = 1408
    5  = 8528
    6  = 64
    7  = 479
    8  = 628
    9  = 247
    10 = 146
    11 = 3430
    13 = 345
    14 = 62
    15 = 979
    18 = 295
    19 = 166
    20 = 2099
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
